$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-321)
# from the old date serial 45172 to the new date serial 45175.
$ws.Range("C2:C321").Value = 45175
